$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 173.4
$ws.Range("I33").Value = 169.75
$ws.Range("K33").Value = 169.75
$ws.Range("M33").Value = 59.25
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H69").Value = 6266.6665
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 6266.6665
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 18799.9995
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -20547.9995
$ws.Range("H70").Value = 1205
$ws.Range("J70").Value = 1205
$ws.Range("L70").Value = 3615
$ws.Range("N70").Value = -4155
$ws.Range("H72").Value = 6266.6665
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 6266.6665
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 56399.9985
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -65135.9985
$ws.Range("H73").Value = 1205
$ws.Range("J73").Value = 1205
$ws.Range("L73").Value = 3615
$ws.Range("N73").Value = -5487
$ws.Range("H115").Value = 425
$ws.Range("I115").Value = 425
$ws.Range("K115").Value = 1275
$ws.Range("M115").Value = 292
$ws.Range("H137").Value = 60623.176
$ws.Range("I137").Value = 1343.1111
$ws.Range("J137").Value = 127313.25
$ws.Range("K137").Value = 4029.3333
$ws.Range("L137").Value = 381939.75
$ws.Range("M137").Value = -1479.3333
$ws.Range("N137").Value = -387039.75
$ws.Range("H138").Value = 3655.8572
$ws.Range("J138").Value = 2960.282
$ws.Range("L138").Value = 8880.846000000001
$ws.Range("N138").Value = -19160.846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 69754.75
$ws.Range("I23").Value = 59670.668
$ws.Range("K23").Value = 59670.668
$ws.Range("M23").Value = -59411.668
$ws.Range("H32").Value = 10997.571
$ws.Range("I32").Value = 8459.833000000001
$ws.Range("J32").Value = 26224
$ws.Range("K32").Value = 8459.833000000001
$ws.Range("L32").Value = 26224
$ws.Range("M32").Value = -8172.833000000001
$ws.Range("N32").Value = -26798
$ws.Range("H37").Value = 12680
$ws.Range("J37").Value = 19000
$ws.Range("L37").Value = 19000
$ws.Range("N37").Value = -19546
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H45").Value = 5295442
$ws.Range("I45").Value = 11250762
$ws.Range("J45").Value = 1823.8889
$ws.Range("K45").Value = 11250762
$ws.Range("L45").Value = 1823.8889
$ws.Range("M45").Value = -11250385
$ws.Range("N45").Value = -2577.8889
$ws.Range("H63").Value = 3099.5
$ws.Range("I63").Value = 2999
$ws.Range("J63").Value = 3200
$ws.Range("K63").Value = 2999
$ws.Range("L63").Value = 3200
$ws.Range("M63").Value = -2313
$ws.Range("N63").Value = -4572
$ws.Range("H66").Value = 3099.5
$ws.Range("I66").Value = 2999
$ws.Range("J66").Value = 3200
$ws.Range("K66").Value = 14995
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = -11563
$ws.Range("N66").Value = -22864
$ws.Range("H97").Value = 921.3333
$ws.Range("I97").Value = 921.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 921.3333
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -425.3333
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 5490.3335
$ws.Range("I26").Value = 5490.3335
$ws.Range("K26").Value = 5490.3335
$ws.Range("M26").Value = -5198.3335
$ws.Range("H94").Value = 248.66667
$ws.Range("I94").Value = 248.66667
$ws.Range("K94").Value = 248.66667
$ws.Range("M94").Value = 202.33333
$ws.Range("H99").Value = 852.5
$ws.Range("I99").Value = 870
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 870
$ws.Range("L99").Value = 800
$ws.Range("M99").Value = 628
$ws.Range("N99").Value = -3796
$ws.Range("H126").Value = 39999
$ws.Range("J126").Value = 39999
$ws.Range("L126").Value = 39999
$ws.Range("N126").Value = -49879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 43000
$ws.Range("J28").Value = 43000
$ws.Range("L28").Value = 43000
$ws.Range("N28").Value = -43490
$ws.Range("H31").Value = 2703.5454
$ws.Range("I31").Value = 2189.45
$ws.Range("J31").Value = 3494.4614
$ws.Range("K31").Value = 2189.45
$ws.Range("L31").Value = 3494.4614
$ws.Range("M31").Value = -1894.45
$ws.Range("N31").Value = -4084.4614
$ws.Range("H34").Value = 2703.5454
$ws.Range("I34").Value = 2189.45
$ws.Range("J34").Value = 3494.4614
$ws.Range("K34").Value = 2189.45
$ws.Range("L34").Value = 3494.4614
$ws.Range("M34").Value = -1987.45
$ws.Range("N34").Value = -3898.4614
$ws.Range("H129").Value = 46249.25
$ws.Range("J129").Value = 46249.25
$ws.Range("L129").Value = 46249.25
$ws.Range("N129").Value = -56249.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 45132.12
$ws.Range("I113").Value = 138398.12
$ws.Range("J113").Value = 1242.2354
$ws.Range("K113").Value = 415194.36
$ws.Range("L113").Value = 3726.7062
$ws.Range("M113").Value = -413024.36
$ws.Range("N113").Value = -8066.706200000001
$ws.Range("H122").Value = 1539.4667
$ws.Range("I122").Value = 550
$ws.Range("J122").Value = 1691.6923
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 15225.2307
$ws.Range("M122").Value = -2500
$ws.Range("N122").Value = -20125.2307
$ws.Range("H136").Value = 1273.8235
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 5809.087
$ws.Range("J137").Value = 7185.8125
$ws.Range("L137").Value = 21557.4375
$ws.Range("N137").Value = -31757.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13033.333
$ws.Range("H70").Value = 5291.857
$ws.Range("J70").Value = 5417
$ws.Range("L70").Value = 5417
$ws.Range("N70").Value = -5957
$ws.Range("H73").Value = 5291.857
$ws.Range("J73").Value = 5417
$ws.Range("L73").Value = 5417
$ws.Range("N73").Value = -7289
$ws.Range("H80").Value = 2850.6
$ws.Range("I80").Value = 2913.1538
$ws.Range("K80").Value = 2913.1538
$ws.Range("M80").Value = -1915.1538
$ws.Range("H83").Value = 2850.6
$ws.Range("I83").Value = 2913.1538
$ws.Range("K83").Value = 14565.769
$ws.Range("M83").Value = -9573.769
$ws.Range("H132").Value = 1044060.94
$ws.Range("I132").Value = 1429607.6
$ws.Range("K132").Value = 4288822.800000001
$ws.Range("M132").Value = -4286292.800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1169.0625
$ws.Range("I22").Value = 847.3333
$ws.Range("J22").Value = 1582.7142
$ws.Range("K22").Value = 847.3333
$ws.Range("L22").Value = 1582.7142
$ws.Range("M22").Value = -552.3333
$ws.Range("N22").Value = -2172.7142
$ws.Range("H27").Value = 1169.0625
$ws.Range("I27").Value = 847.3333
$ws.Range("J27").Value = 1582.7142
$ws.Range("K27").Value = 847.3333
$ws.Range("L27").Value = 1582.7142
$ws.Range("M27").Value = -740.3333
$ws.Range("N27").Value = -1796.7142
$ws.Range("H40").Value = 23662.75
$ws.Range("J40").Value = 13459.6
$ws.Range("L40").Value = 13459.6
$ws.Range("N40").Value = -13731.6
$ws.Range("H46").Value = 2013.6923
$ws.Range("I46").Value = 1798
$ws.Range("J46").Value = 2148.5
$ws.Range("K46").Value = 1798
$ws.Range("L46").Value = 2148.5
$ws.Range("M46").Value = -1610
$ws.Range("N46").Value = -2524.5
$ws.Range("H68").Value = 2476.9333
$ws.Range("I68").Value = 1923.091
$ws.Range("K68").Value = 1923.091
$ws.Range("M68").Value = -1174.091
$ws.Range("H71").Value = 2476.9333
$ws.Range("I71").Value = 1923.091
$ws.Range("K71").Value = 9615.455
$ws.Range("M71").Value = -5871.455
$ws.Range("H82").Value = 1664.8125
$ws.Range("I82").Value = 1731.3334
$ws.Range("J82").Value = 1465.25
$ws.Range("K82").Value = 1731.3334
$ws.Range("L82").Value = 1465.25
$ws.Range("M82").Value = -1370.3334
$ws.Range("N82").Value = -2187.25
$ws.Range("H85").Value = 1664.8125
$ws.Range("I85").Value = 1731.3334
$ws.Range("J85").Value = 1465.25
$ws.Range("K85").Value = 1731.3334
$ws.Range("L85").Value = 1465.25
$ws.Range("M85").Value = -483.3334
$ws.Range("N85").Value = -3961.25
$ws.Range("H100").Value = 1858.6
$ws.Range("I100").Value = 1623.25
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 1623.25
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -1082.25
$ws.Range("N100").Value = -3882

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 774.5
$ws.Range("I81").Value = 774.5
$ws.Range("J81").Value = 774.5
$ws.Range("K81").Value = 1549
$ws.Range("L81").Value = 1549
$ws.Range("M81").Value = -488
$ws.Range("N81").Value = -3671
$ws.Range("H84").Value = 774.5
$ws.Range("I84").Value = 774.5
$ws.Range("J84").Value = 774.5
$ws.Range("K84").Value = 7745
$ws.Range("L84").Value = 7745
$ws.Range("M84").Value = -2441
$ws.Range("N84").Value = -18353
$ws.Range("H136").Value = 14621650
$ws.Range("I136").Value = 21368906
$ws.Range("K136").Value = 64106718
$ws.Range("M136").Value = -64104168
